$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2774.34
$ws.Range("I64").Value = 2560.9033
$ws.Range("K64").Value = 2560.9033
$ws.Range("M64").Value = -2312.9033
$ws.Range("H67").Value = 2774.34
$ws.Range("I67").Value = 2560.9033
$ws.Range("K67").Value = 2560.9033
$ws.Range("M67").Value = -1702.9033
$ws.Range("H121").Value = 6160.2
$ws.Range("J121").Value = 7625
$ws.Range("L121").Value = 22875
$ws.Range("N121").Value = -26369
$ws.Range("H141").Value = 1184.963
$ws.Range("I141").Value = 879.76
$ws.Range("K141").Value = 2639.28
$ws.Range("M141").Value = 2540.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1744.32
$ws.Range("I2").Value = 1555.4
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 1555.4
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = -1442.4
$ws.Range("N2").Value = -2726
$ws.Range("H32").Value = 4703.3022
$ws.Range("I32").Value = 2791.4321
$ws.Range("J32").Value = 15027.4
$ws.Range("K32").Value = 2791.4321
$ws.Range("L32").Value = 15027.4
$ws.Range("M32").Value = -2504.4321
$ws.Range("N32").Value = -15601.4
$ws.Range("H74").Value = 28596.703
$ws.Range("I74").Value = 34492.535
$ws.Range("K74").Value = 34492.535
$ws.Range("M74").Value = -33618.535
$ws.Range("H77").Value = 28596.703
$ws.Range("I77").Value = 34492.535
$ws.Range("K77").Value = 172462.675
$ws.Range("M77").Value = -168094.675
$ws.Range("H101").Value = 26300.2
$ws.Range("J101").Value = 26300.2
$ws.Range("L101").Value = 26300.2
$ws.Range("N101").Value = -32790.2
$ws.Range("H102").Value = 1973.6666
$ws.Range("I102").Value = 1855
$ws.Range("J102").Value = 2211
$ws.Range("K102").Value = 1855
$ws.Range("L102").Value = 2211
$ws.Range("M102").Value = -233
$ws.Range("N102").Value = -5455
$ws.Range("H116").Value = 1744.32
$ws.Range("I116").Value = 1555.4
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 1555.4
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = 738.5999999999999
$ws.Range("N116").Value = -7088
$ws.Range("H132").Value = 2216
$ws.Range("I132").Value = 1805.2
$ws.Range("J132").Value = 2900.6667
$ws.Range("K132").Value = 5415.6
$ws.Range("L132").Value = 8702.000100000001
$ws.Range("M132").Value = -2885.6
$ws.Range("N132").Value = -13762.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1744.32
$ws.Range("I3").Value = 1555.4
$ws.Range("J3").Value = 2500
$ws.Range("K3").Value = 1555.4
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = -1441.4
$ws.Range("N3").Value = -2728
$ws.Range("H100").Value = 22877.572
$ws.Range("J100").Value = 22877.572
$ws.Range("L100").Value = 22877.572
$ws.Range("N100").Value = -25041.572
$ws.Range("H107").Value = 1167.4
$ws.Range("I107").Value = 1074
$ws.Range("J107").Value = 1696.6666
$ws.Range("K107").Value = 1074
$ws.Range("L107").Value = 1696.6666
$ws.Range("M107").Value = 846
$ws.Range("N107").Value = -5536.6666
$ws.Range("H134").Value = 693170.5
$ws.Range("I134").Value = 1003157.5
$ws.Range("K134").Value = 3009472.5
$ws.Range("M134").Value = -3006937.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10205961
$ws.Range("I31").Value = 1268.069
$ws.Range("J31").Value = 25002766
$ws.Range("K31").Value = 1268.069
$ws.Range("L31").Value = 25002766
$ws.Range("M31").Value = -973.069
$ws.Range("N31").Value = -25003356
$ws.Range("H34").Value = 10205961
$ws.Range("I34").Value = 1268.069
$ws.Range("J34").Value = 25002766
$ws.Range("K34").Value = 1268.069
$ws.Range("L34").Value = 25002766
$ws.Range("M34").Value = -1066.069
$ws.Range("N34").Value = -25003170
$ws.Range("H134").Value = 2569.9443
$ws.Range("I134").Value = 2781.577
$ws.Range("J134").Value = 2019.7
$ws.Range("K134").Value = 8344.731
$ws.Range("L134").Value = 6059.1
$ws.Range("M134").Value = -5809.731
$ws.Range("N134").Value = -11129.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 971.45
$ws.Range("I97").Value = 452.94116
$ws.Range("J97").Value = 1354.6957
$ws.Range("K97").Value = 1358.82348
$ws.Range("L97").Value = 4064.0871
$ws.Range("M97").Value = -862.82348
$ws.Range("N97").Value = -5056.0871
$ws.Range("H98").Value = 542.1429000000001
$ws.Range("I98").Value = 493.57144
$ws.Range("J98").Value = 590.7143
$ws.Range("K98").Value = 1480.71432
$ws.Range("L98").Value = 1772.1429
$ws.Range("M98").Value = 17.28567999999996
$ws.Range("N98").Value = -4768.1429
$ws.Range("H113").Value = 3367527.5
$ws.Range("I113").Value = 5051043.5
$ws.Range("J113").Value = 496
$ws.Range("K113").Value = 15153130.5
$ws.Range("L113").Value = 1488
$ws.Range("M113").Value = -15150960.5
$ws.Range("N113").Value = -5828
$ws.Range("H134").Value = 8976.076999999999
$ws.Range("I134").Value = 8384.143
$ws.Range("J134").Value = 9666.666999999999
$ws.Range("K134").Value = 25152.429
$ws.Range("L134").Value = 29000.001
$ws.Range("M134").Value = -20082.429
$ws.Range("N134").Value = -39140.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 28277.5
$ws.Range("J64").Value = 28277.5
$ws.Range("L64").Value = 28277.5
$ws.Range("N64").Value = -28773.5
$ws.Range("H67").Value = 28277.5
$ws.Range("J67").Value = 28277.5
$ws.Range("L67").Value = 28277.5
$ws.Range("N67").Value = -29993.5
$ws.Range("H113").Value = 1471.4286
$ws.Range("I113").Value = 1377.7778
$ws.Range("J113").Value = 1640
$ws.Range("K113").Value = 1377.7778
$ws.Range("L113").Value = 1640
$ws.Range("M113").Value = 792.2221999999999
$ws.Range("N113").Value = -5980

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2480.8
$ws.Range("I40").Value = 2476
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 2476
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -2340
$ws.Range("N40").Value = -2772
$ws.Range("H46").Value = 879
$ws.Range("I46").Value = 967.3333
$ws.Range("J46").Value = 829.3125
$ws.Range("K46").Value = 967.3333
$ws.Range("L46").Value = 829.3125
$ws.Range("M46").Value = -779.3333
$ws.Range("N46").Value = -1205.3125
$ws.Range("H132").Value = 4127.9375
$ws.Range("I132").Value = 4016.1562
$ws.Range("J132").Value = 4351.5
$ws.Range("K132").Value = 12048.4686
$ws.Range("L132").Value = 13054.5
$ws.Range("M132").Value = -9518.4686
$ws.Range("N132").Value = -18114.5
$ws.Range("H133").Value = 31315
$ws.Range("J133").Value = 31315
$ws.Range("L133").Value = 31315
$ws.Range("N133").Value = -36375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 44933.332
$ws.Range("J97").Value = 44933.332
$ws.Range("L97").Value = 44933.332
$ws.Range("N97").Value = -46915.332
$ws.Range("H101").Value = 15073.75
$ws.Range("J101").Value = 15073.75
$ws.Range("L101").Value = 15073.75
$ws.Range("N101").Value = -21563.75
$ws.Range("H113").Value = 746.4
$ws.Range("I113").Value = 593.25
$ws.Range("J113").Value = 921.4286
$ws.Range("K113").Value = 1779.75
$ws.Range("L113").Value = 2764.2858
$ws.Range("M113").Value = 390.25
$ws.Range("N113").Value = -7104.2858
$ws.Range("H136").Value = 1804.081
$ws.Range("I136").Value = 1522.4182
$ws.Range("J136").Value = 2619.4211
$ws.Range("K136").Value = 4567.2546
$ws.Range("L136").Value = 7858.263300000001
$ws.Range("M136").Value = -2017.2546
$ws.Range("N136").Value = -12958.2633
